# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the
# b76cd5e1-fd44-4e10-b780-c1cfe23e7892 row (row 4) on both the
# zh-cn and de-de sheets to reflect the newly generated report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 22:37:59"
$wsZhCn.Range("H4").Value = "2016-03-12 22:38:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 22:38:02"
$wsDeDe.Range("H4").Value = "2016-03-12 22:38:23"
